# "Ran Priors with _1"
# Update the prior variance (column D) and standard deviation/weight (column E)
# for the data rows (2 through 44) of the Priors table, and update the
# view/selection to reflect where the user ended up scrolled to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-44: column D (variance multiplier) changes from 1 to 5,
# and column E (prior std/weight) is unified to 0.2 (was 0.15 for rows 2-22
# and 0.05 for rows 23-44).
$ws.Range("D2:D44").Value = 5
$ws.Range("E2:E44").Value = 0.2

# Reflect the scrolled viewport / active cell selection recorded after the run.
$ws.Application.Goto($ws.Range("A22"), $true) | Out-Null
$ws.Range("E33").Select() | Out-Null
